# Update gh-pages to output generated at 456a3b4
# This script updates the "F" column (想去人数 / want-to-go count) values
# across the four worksheets to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 419
$ws.Range("F4").Value  = 1178
$ws.Range("F5").Value  = 60
$ws.Range("F7").Value  = 44
$ws.Range("F8").Value  = 1082
$ws.Range("F10").Value = 391
$ws.Range("F11").Value = 438
$ws.Range("F13").Value = 324
$ws.Range("F14").Value = 371
$ws.Range("F15").Value = 53
$ws.Range("F18").Value = 586
$ws.Range("F19").Value = 1481
$ws.Range("F20").Value = 5776
$ws.Range("F22").Value = 1628
$ws.Range("F23").Value = 388
$ws.Range("F24").Value = 76
$ws.Range("F25").Value = 36
$ws.Range("F26").Value = 5418
$ws.Range("F27").Value = 5418
$ws.Range("F28").Value = 134
$ws.Range("F30").Value = 1558
$ws.Range("F31").Value = 188
$ws.Range("F32").Value = 28
$ws.Range("F33").Value = 69
$ws.Range("F34").Value = 1082
$ws.Range("F36").Value = 125
$ws.Range("F37").Value = 8
$ws.Range("F38").Value = 71
$ws.Range("F39").Value = 3818

# ---------------------------------------------------------------------------
# Sheet: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 18
$ws.Range("F4").Value = 68
$ws.Range("F5").Value = 176
$ws.Range("F8").Value = 249

# ---------------------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9430
$ws.Range("F4").Value = 2166
$ws.Range("F5").Value = 238

# ---------------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 9430
$ws.Range("F4").Value  = 2166
$ws.Range("F6").Value  = 419
$ws.Range("F7").Value  = 1178
$ws.Range("F8").Value  = 60
$ws.Range("F10").Value = 44
$ws.Range("F11").Value = 1082
$ws.Range("F12").Value = 391
$ws.Range("F13").Value = 438
$ws.Range("F14").Value = 324
$ws.Range("F15").Value = 371
$ws.Range("F16").Value = 53
$ws.Range("F19").Value = 18
$ws.Range("F20").Value = 1481
$ws.Range("F21").Value = 5776
$ws.Range("F23").Value = 1628
$ws.Range("F26").Value = 388
$ws.Range("F29").Value = 5418
$ws.Range("F30").Value = 5418
$ws.Range("F31").Value = 134
$ws.Range("F33").Value = 1558
$ws.Range("F34").Value = 188
$ws.Range("F35").Value = 28
$ws.Range("F36").Value = 1082
$ws.Range("F38").Value = 125
$ws.Range("F44").Value = 71
$ws.Range("F46").Value = 3818
